$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 218.7619
$ws.Range("I33").Value = 230
$ws.Range("K33").Value = 230
$ws.Range("M33").Value = -1
$ws.Range("H52").Value = 1494
$ws.Range("I52").Value = 1772.3334
$ws.Range("K52").Value = 5317.0002
$ws.Range("M52").Value = -5157.0002
$ws.Range("H62").Value = 29635144
$ws.Range("I62").Value = 33338288
$ws.Range("K62").Value = 33338288
$ws.Range("M62").Value = -33337664
$ws.Range("H65").Value = 29635144
$ws.Range("I65").Value = 33338288
$ws.Range("K65").Value = 166691440
$ws.Range("M65").Value = -166688320
$ws.Range("H75").Value = 148230.08
$ws.Range("I75").Value = 192207
$ws.Range("J75").Value = 120744.5
$ws.Range("K75").Value = 192207
$ws.Range("L75").Value = 120744.5
$ws.Range("M75").Value = -191271
$ws.Range("N75").Value = -122616.5
$ws.Range("H78").Value = 148230.08
$ws.Range("I78").Value = 192207
$ws.Range("J78").Value = 120744.5
$ws.Range("K78").Value = 576621
$ws.Range("L78").Value = 362233.5
$ws.Range("M78").Value = -571941
$ws.Range("N78").Value = -371593.5
$ws.Range("H113").Value = 52539.35
$ws.Range("J113").Value = 127300.625
$ws.Range("L113").Value = 127300.625
$ws.Range("N113").Value = -133808.625
$ws.Range("H116").Value = 5249.1055
$ws.Range("J116").Value = 6763.5
$ws.Range("L116").Value = 6763.5
$ws.Range("N116").Value = -13647.5
$ws.Range("H120").Value = 146666.33
$ws.Range("J120").Value = 146666.33
$ws.Range("L120").Value = 146666.33
$ws.Range("N120").Value = -156342.33
$ws.Range("H132").Value = 995.4194
$ws.Range("I132").Value = 933.0345
$ws.Range("K132").Value = 2799.1035
$ws.Range("M132").Value = -269.1035000000002
$ws.Range("H137").Value = 3288.8696
$ws.Range("I137").Value = 2578.4517
$ws.Range("K137").Value = 7735.355100000001
$ws.Range("M137").Value = -5185.355100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 42234.4
$ws.Range("I32").Value = 45241.965
$ws.Range("J32").Value = 15166.333
$ws.Range("K32").Value = 45241.965
$ws.Range("L32").Value = 15166.333
$ws.Range("M32").Value = -44954.965
$ws.Range("N32").Value = -15740.333
$ws.Range("H74").Value = 402677.2
$ws.Range("I74").Value = 834687.5600000001
$ws.Range("J74").Value = 3898.3845
$ws.Range("K74").Value = 834687.5600000001
$ws.Range("L74").Value = 3898.3845
$ws.Range("M74").Value = -833813.5600000001
$ws.Range("N74").Value = -5646.3845
$ws.Range("H77").Value = 402677.2
$ws.Range("I77").Value = 834687.5600000001
$ws.Range("J77").Value = 3898.3845
$ws.Range("K77").Value = 4173437.8
$ws.Range("L77").Value = 19491.9225
$ws.Range("M77").Value = -4169069.8
$ws.Range("N77").Value = -28227.9225
$ws.Range("H81").Value = 50000
$ws.Range("J81").Value = 50000
$ws.Range("L81").Value = 50000
$ws.Range("N81").Value = -51996
$ws.Range("H84").Value = 50000
$ws.Range("J84").Value = 50000
$ws.Range("L84").Value = 150000
$ws.Range("N84").Value = -159984
$ws.Range("H110").Value = 20835012
$ws.Range("J110").Value = 2561.8572
$ws.Range("L110").Value = 2561.8572
$ws.Range("N110").Value = -6651.8572

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 6075.2563
$ws.Range("I134").Value = 4394.4585
$ws.Range("J134").Value = 8764.532999999999
$ws.Range("K134").Value = 13183.3755
$ws.Range("L134").Value = 26293.599
$ws.Range("M134").Value = -10648.3755
$ws.Range("N134").Value = -31363.599

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5810.7915
$ws.Range("I31").Value = 3681
$ws.Range("K31").Value = 3681
$ws.Range("M31").Value = -3386
$ws.Range("H34").Value = 5810.7915
$ws.Range("I34").Value = 3681
$ws.Range("K34").Value = 3681
$ws.Range("M34").Value = -3479
$ws.Range("H50").Value = 42498.25
$ws.Range("J50").Value = 42498.25
$ws.Range("L50").Value = 42498.25
$ws.Range("N50").Value = -43748.25
$ws.Range("H132").Value = 29671.348
$ws.Range("I132").Value = 3271.6924
$ws.Range("K132").Value = 9815.0772
$ws.Range("M132").Value = -7285.0772
$ws.Range("H134").Value = 8540.056
$ws.Range("I134").Value = 9137
$ws.Range("J134").Value = 7793.875
$ws.Range("K134").Value = 27411
$ws.Range("L134").Value = 23381.625
$ws.Range("M134").Value = -24876
$ws.Range("N134").Value = -28451.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 1534
$ws.Range("I11").Value = 40.333332
$ws.Range("J11").Value = 1982.1
$ws.Range("K11").Value = 120.999996
$ws.Range("L11").Value = 5946.299999999999
$ws.Range("M11").Value = 19.000004
$ws.Range("N11").Value = -6226.299999999999
$ws.Range("H19").Value = 492
$ws.Range("I19").Value = 100
$ws.Range("K19").Value = 300
$ws.Range("M19").Value = -126
$ws.Range("H22").Value = 675
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
$ws.Range("H27").Value = 675
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
$ws.Range("H37").Value = 99810.28999999999
$ws.Range("J37").Value = 99810.28999999999
$ws.Range("L37").Value = 299430.87
$ws.Range("N37").Value = -299654.87
$ws.Range("H46").Value = 175
$ws.Range("I46").Value = 100
$ws.Range("K46").Value = 300
$ws.Range("M46").Value = -209
$ws.Range("H92").Value = 657.7
$ws.Range("I92").Value = 711.8182
$ws.Range("J92").Value = 591.55554
$ws.Range("K92").Value = 2135.4546
$ws.Range("L92").Value = 1774.66662
$ws.Range("M92").Value = -887.4546
$ws.Range("N92").Value = -4270.66662
$ws.Range("H113").Value = 2794.0715
$ws.Range("I113").Value = 3568.5
$ws.Range("J113").Value = 2484.3
$ws.Range("K113").Value = 10705.5
$ws.Range("L113").Value = 7452.900000000001
$ws.Range("M113").Value = -8535.5
$ws.Range("N113").Value = -11792.9
$ws.Range("H134").Value = 2728.9285
$ws.Range("J134").Value = 8022
$ws.Range("L134").Value = 24066
$ws.Range("N134").Value = -34206
$ws.Range("H137").Value = 8357.1875
$ws.Range("I137").Value = 9849.909
$ws.Range("J137").Value = 5073.2
$ws.Range("K137").Value = 29549.727
$ws.Range("L137").Value = 15219.6
$ws.Range("M137").Value = -24449.727
$ws.Range("N137").Value = -25419.6
$ws.Range("H140").Value = 1761.9333
$ws.Range("J140").Value = 2042.3636
$ws.Range("L140").Value = 6127.0908
$ws.Range("N140").Value = -16487.0908

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 639.4091
$ws.Range("I16").Value = 556.2105
$ws.Range("J16").Value = 1166.3334
$ws.Range("K16").Value = 556.2105
$ws.Range("L16").Value = 1166.3334
$ws.Range("M16").Value = -386.2105
$ws.Range("N16").Value = -1506.3334
$ws.Range("H22").Value = 3462.1462
$ws.Range("I22").Value = 2548.0527
$ws.Range("J22").Value = 4251.591
$ws.Range("K22").Value = 2548.0527
$ws.Range("L22").Value = 4251.591
$ws.Range("M22").Value = -2253.0527
$ws.Range("N22").Value = -4841.591
$ws.Range("H27").Value = 3462.1462
$ws.Range("I27").Value = 2548.0527
$ws.Range("J27").Value = 4251.591
$ws.Range("K27").Value = 2548.0527
$ws.Range("L27").Value = 4251.591
$ws.Range("M27").Value = -2441.0527
$ws.Range("N27").Value = -4465.591
$ws.Range("H55").Value = 175.03847
$ws.Range("I55").Value = 184.25
$ws.Range("J55").Value = 167.14285
$ws.Range("K55").Value = 184.25
$ws.Range("L55").Value = 167.14285
$ws.Range("M55").Value = -11.25
$ws.Range("N55").Value = -513.14285
$ws.Range("H61").Value = 3380.2
$ws.Range("I61").Value = 2725.25
$ws.Range("K61").Value = 2725.25
$ws.Range("M61").Value = -2523.25
$ws.Range("H113").Value = 3380.2
$ws.Range("I113").Value = 2725.25
$ws.Range("K113").Value = 2725.25
$ws.Range("M113").Value = -555.25
$ws.Range("H132").Value = 4836.914
$ws.Range("I132").Value = 3431.1052
$ws.Range("K132").Value = 10293.3156
$ws.Range("M132").Value = -7763.3156
$ws.Range("H133").Value = 111111
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 111111
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 111111
$ws.Range("N133").Value = -116171
$ws.Range("M133").ClearContents()
$ws.Range("H136").Value = 3037948.5
$ws.Range("I136").Value = 3515893
$ws.Range("K136").Value = 10547679
$ws.Range("M136").Value = -10545129

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 13621.3
$ws.Range("J69").Value = 13621.3
$ws.Range("L69").Value = 13621.3
$ws.Range("N69").Value = -15119.3
$ws.Range("H72").Value = 13621.3
$ws.Range("J72").Value = 13621.3
$ws.Range("L72").Value = 40863.89999999999
$ws.Range("N72").Value = -48351.89999999999
$ws.Range("H80").Value = 25500
$ws.Range("J80").Value = 25500
$ws.Range("L80").Value = 25500
$ws.Range("N80").Value = -27496
$ws.Range("H83").Value = 25500
$ws.Range("J83").Value = 25500
$ws.Range("L83").Value = 76500
$ws.Range("N83").Value = -86484
$ws.Range("H122").Value = 3372.6428
$ws.Range("I122").Value = 3169.7222
$ws.Range("K122").Value = 9509.1666
$ws.Range("M122").Value = -7059.1666
$ws.Range("H136").Value = 3099.3572
$ws.Range("I136").Value = 1626.8334
$ws.Range("J136").Value = 5749.9
$ws.Range("K136").Value = 4880.5002
$ws.Range("L136").Value = 17249.7
$ws.Range("M136").Value = -2330.5002
$ws.Range("N136").Value = -22349.7
